$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
# ---------------------------------------------------------------------
$titlePar = $d.Paragraphs.Item(1)
$titlePar.Range.InsertParagraphAfter()

$metaPar = $d.Paragraphs.Item(2)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r/>' +
           '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
           '<w:r><w:t>: Discover the legend of King Arthur with Avalon. Enjoy bonuses, free spins, and Wilds with this online slot machine. Play for free now.</w:t></w:r>' +
           '</w:p>'
$metaPar.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicate bold title
#    paragraph and replace the italic meta-description paragraph's
#    text with the new image-generation prompt (keep its formatting).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPar = $d.Paragraphs.Item($count - 1)
$italicPar = $d.Paragraphs.Item($count)

# Delete the whole bold paragraph (including its paragraph mark).
$dupRange = $d.Range($boldPar.Range.Start, $italicPar.Range.Start)
$dupRange.Delete()

# Replace the text of the (now last) italic paragraph, keeping its
# paragraph mark / run formatting untouched.
$lastPar = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPrompt = 'Please create a cartoon-style feature image of a happy Maya warrior with glasses for the online slot game "Avalon". The image should be bright and eye-catching, with the character holding a sword or treasure chest to represent the game''s theme of King Arthur''s legend and Island. You can add elements like symbols, crowns, and poker cards to the image to showcase the game''s graphics and symbols. The Maya warrior should be smiling and look happy, adding a playful and engaging tone to the image. The background should highlight the island setting, with forests, mountains, and an aura of mystery and excitement. Overall, the image should be visually appealing and capture the essence of the game''s theme and gameplay.'

$textRange = $d.Range($lastPar.Range.Start, $lastPar.Range.End - 1)
$textRange.Text = $newPrompt

Write-Output "Done"
